$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.011.78'
$ws.Range('E2').Value = '  +0.30%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.790.01'
$ws.Range('E3').Value = '  -0.42%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.31%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '600.67'
$ws.Range('E5').Value = '  +0.63%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '165.37'
$ws.Range('E6').Value = '  -1.01%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.517'
$ws.Range('E8').Value = '  -0.79%  '

$ws.Range('E9').Value = '  -1.08%  '

$ws.Range('E10').Value = '  +0.32%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.46'
$ws.Range('E11').Value = '  +2.66%  '

$ws.Range('E12').Value = '  -1.67%  '

$ws.Range('E13').Value = '  -1.02%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.424.05'
$ws.Range('E14').Value = '  -0.44%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.783.00'
$ws.Range('E15').Value = '  -1.53%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.997.74'
$ws.Range('E16').Value = '  +0.30%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '18.40'
$ws.Range('E17').Value = '  -1.30%  '

$ws.Range('E18').Value = '  +2.00%  '

$ws.Range('E19').Value = '  -0.65%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '460.66'
$ws.Range('E20').Value = '  -0.18%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.70'
$ws.Range('E21').Value = '  -1.77%  '

$ws.Range('E22').Value = '  -0.56%  '

$ws.Range('E23').Value = '  -2.94%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.78'
$ws.Range('E24').Value = '  -0.88%  '

$ws.Range('E25').Value = '  -0.15%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.10'
$ws.Range('E26').Value = '  +0.17%  '

$ws.Range('E27').Value = '  -0.17%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.98'
$ws.Range('E28').Value = '  -0.41%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.937.05'
$ws.Range('E29').Value = '  -0.36%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.38'
$ws.Range('E30').Value = '  +1.85%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.63'
$ws.Range('E31').Value = '  -5.68%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '29.30'
$ws.Range('E33').Value = '  -1.41%  '

$ws.Range('E34').Value = '  +0.02%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '8.98'
$ws.Range('E35').Value = '  -1.08%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.100'
$ws.Range('E36').Value = '  -0.01%  '

$ws.Range('E37').Value = '  +0.35%  '

$ws.Range('E38').Value = '  -3.41%  '

$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.989'
$ws.Range('E39').Value = '  -0.54%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.77'
$ws.Range('E40').Value = '  -0.27%  '

$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('E43').Value = '  +0.42%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '47.36'
$ws.Range('E44').Value = '  -1.58%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '42.98'
$ws.Range('E45').Value = '  -2.16%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '151.56'
$ws.Range('E46').Value = '  +0.48%  '

$ws.Range('E47').Value = '  +0.31%  '

$ws.Range('E48').Value = '  +2.46%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '391.99'
$ws.Range('E49').Value = '  -0.13%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.34'
$ws.Range('E50').Value = '  +5.90%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '26.76'
$ws.Range('E51').Value = '  +1.33%  '
